$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 0. The document currently ends with a hidden "_GoBack" bookmark (right
#    after the inline picture). That bookmark is being relocated to the new
#    bullet point we are about to add, so remove the old one first. Doing
#    this before inserting the new bookmark avoids any ambiguity from two
#    bookmarks sharing the same name at the same time.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 1. Italicize "ubuntu-test-server" in the existing "[targets]" block.
# ---------------------------------------------------------------------------
$targetsLine = $d.Paragraphs(7).Range
$targetsLine.Find.Execute("ubuntu-test-server", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetsLine.Font.Italic = $true

# ---------------------------------------------------------------------------
# 2. Italicize "mkeith" on the same line.
# ---------------------------------------------------------------------------
$mkeithRun = $d.Paragraphs(7).Range
$mkeithRun.Find.Execute("mkeith", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mkeithRun.Font.Italic = $true

# ---------------------------------------------------------------------------
# 3. Insert a brand-new bulleted paragraph right after the blank "NoSpacing"
#    paragraph that follows the "...ansible_user=mkeith" line, i.e. right
#    before the "Run the script: ..." bullet.
# ---------------------------------------------------------------------------
$blankAfterMkeith = $d.Paragraphs(8)
$insertionPoint = $blankAfterMkeith.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$newParagraph = $d.Paragraphs(9)

$newParagraphXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Be sure to change the hosts in the .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yml</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file as we </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">hosts: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>ubuntu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t>-test-server</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newParagraph.Range.InsertXML($newParagraphXml)
